# "Generate Report for Handoff"
# The localization-status report was regenerated: the previous handback
# status/timestamps are replaced with a fresh "Ready for handoff" status
# and updated generation timestamps, and the (now shorter) status columns
# are re-sized to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed generation / handoff timestamps
$overview.Range("G2").Value = "2016-09-08 05:18:21"
$dede.Range("H2").Value     = "2016-09-08 05:18:21"
$zhcn.Range("H2").Value     = "2016-09-08 05:18:16"

# --- Column widths re-fit to the new (shorter) status text.
# ColumnWidth is expressed in characters and snaps to the host's pixel grid,
# so we pick the character width whose rounded result lands on the nearest
# achievable value to the refreshed report's column width.
$overview.Range("E1").ColumnWidth = 16.333333333333332
$overview.Range("F1").ColumnWidth = 16.333333333333332
$zhcn.Range("C1").ColumnWidth     = 29.166666666666668
$dede.Range("C1").ColumnWidth     = 29.166666666666668
